# Update public EPEX spot / Gaz / CO2 sheets with the new day (29-jun / 2025-06-27).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add column P ("29-jun") with header + 24 hourly values.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (O1) onto the new header
# cell (P1) so it keeps the same bold/border/centered style, then set value.
$wsSpot.Range("O1").Copy()
$wsSpot.Range("P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSpot.Range("P1").Value = "29-jun"

$spotValues = @{
    2  = 95.06999999999999
    3  = 89.42
    4  = 89.66
    5  = 86.81999999999999
    6  = 84.3
    7  = 85.3
    8  = 86
    9  = 74.70999999999999
    10 = 47.4
    11 = 3.75
    12 = 0
    13 = -0.08
    14 = -0.02
    15 = -0.03
    16 = -0.09
    17 = -0.01
    18 = 0
    19 = 20
    20 = 79.38
    21 = 102.14
    22 = 111.98
    23 = 118.07
    24 = 113.73
    25 = 101
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 16).Value = $spotValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 13 for 2025-06-27.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to be entered as literal text (matching the existing
# "Date" column cells, which are plain strings, not real dates) by briefly
# switching the destination cell to a text format before assignment, then
# copy the neighboring cell's (A12) formatting back over it so no stray
# date/text number-format sticks around on the new cell.
$wsGaz.Range("A13").NumberFormat = "@"
$wsGaz.Range("A13").Value = "2025-06-27"
$wsGaz.Range("A12").Copy()
$wsGaz.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsGaz.Range("B13").Value = 32.7

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 13 for 2025-06-27.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A13").NumberFormat = "@"
$wsCo2.Range("A13").Value = "2025-06-27"
$wsCo2.Range("A12").Copy()
$wsCo2.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsCo2.Range("B13").Value = 69.92
